$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 7; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 35; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 50; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 58; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 61; I = '%'; J = 'Uninterpretable' },
    @{ Row = 71; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 75; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 77; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 83; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 87; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 101; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 104; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 109; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 115; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 130; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 138; I = 'ba'; J = 'Appreciation' },
    @{ Row = 146; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 147; I = '%'; J = 'Uninterpretable' },
    @{ Row = 149; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 150; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 161; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 163; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 165; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 166; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 168; I = 'ba'; J = 'Appreciation' },
    @{ Row = 179; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 229; I = 'ba'; J = 'Appreciation' },
    @{ Row = 236; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 243; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 248; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 264; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 276; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 277; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 287; I = '%'; J = 'Uninterpretable' },
    @{ Row = 289; I = '%'; J = 'Uninterpretable' },
    @{ Row = 294; I = 'qy'; J = 'Yes-No-Question' },
    @{ Row = 301; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 328; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 349; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 352; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 354; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 365; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 372; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 374; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 393; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 396; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 409; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 411; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 420; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 427; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 429; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 431; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 450; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 457; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 468; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 482; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 499; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 515; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 516; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 566; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 573; I = 'ba'; J = 'Appreciation' },
    @{ Row = 592; I = '%'; J = 'Uninterpretable' },
    @{ Row = 594; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 611; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 617; I = 'ba'; J = 'Appreciation' },
    @{ Row = 630; I = 'ba'; J = 'Appreciation' },
    @{ Row = 655; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 658; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 665; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 9).Value = $u.I
    $ws.Cells.Item($r, 10).Value = $u.J
}
